# Insert a new data row before existing row 33 (shifts rows 33:88 down to 34:89)
# and populate the new row 33 with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("33:33").Insert()

$ws.Range("A33").Value = 7
$ws.Range("B33").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C33").Value = 'Ñuble'
$ws.Range("D33").Value = 44775
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112031
$ws.Range("G33").Value = 'Poroto verde'
$ws.Range("H33").Value = 'Magnum'
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 30000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = 30000
$ws.Range("N33").Value = '$/malla 25 kilos'
$ws.Range("O33").Value = 'Perú'
$ws.Range("P33").Value = 1200
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = 'Hortaliza'
